# Update BunkerPrices at 2025-03-28 06:08
# - Change the number format of S13 (existing "Date" column cell) from a
#   date-only format to a date-time format (style s="3" -> s="2").
# - Append a new data row (row 14) with bunker price values for 2025-03-27,
#   extending the used range from A1:AV13 to A1:AV14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) S13 switches from the "YYYY-MM-DD" format to "YYYY-MM-DD HH:MM:SS" format.
$ws.Range("S13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2) Append new row 14 with the new bunker price data (columns A..AV).
$rowValues = @(521, 572, 536, 522, 652, 530, 629, 665, 565, 526, 568, 519, 580, 527, 651, 772, 560, 668, 45743, 592, 600, 610, 515, 519, 550, 762, 542, 573, 526, 653, 606.5, 579, 533, 576, 884, 650, 510, 629, 546, 521, 535, 519, 515, 497, 530, 552, 497, 561)

$arr = New-Object 'object[,]' 1,$rowValues.Length
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $arr[0, $i] = $rowValues[$i]
}

$ws.Range("A14:AV14").Value = $arr

# The "Date" column (S) on the new row keeps the original date-only format.
$ws.Range("S14").NumberFormat = "YYYY-MM-DD"
